# Scheduled market-price refresh for the Leve profit sheets (ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR).
# Updates currentAveragePrice* / LevePrice* / LeveProfit* columns (H:N) with the latest Universalis pull.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 8463.333000000001
$ws.Range("I51").Value = 5850
$ws.Range("J51").Value = 9770
$ws.Range("K51").Value = 5850
$ws.Range("L51").Value = 9770
$ws.Range("M51").Value = -5366
$ws.Range("N51").Value = -10738
$ws.Range("H62").Value = 28405.516
$ws.Range("I62").Value = 48424.125
$ws.Range("J62").Value = 9564.471
$ws.Range("K62").Value = 48424.125
$ws.Range("L62").Value = 9564.471
$ws.Range("M62").Value = -47800.125
$ws.Range("N62").Value = -10812.471
$ws.Range("H65").Value = 28405.516
$ws.Range("I65").Value = 48424.125
$ws.Range("J65").Value = 9564.471
$ws.Range("K65").Value = 242120.625
$ws.Range("L65").Value = 47822.355
$ws.Range("M65").Value = -239000.625
$ws.Range("N65").Value = -54062.355
$ws.Range("H70").Value = 1600
$ws.Range("I70").Value = 1600
$ws.Range("J70").Value = 1600
$ws.Range("K70").Value = 4800
$ws.Range("L70").Value = 4800
$ws.Range("M70").Value = -4530
$ws.Range("N70").Value = -5340
$ws.Range("H73").Value = 1600
$ws.Range("I73").Value = 1600
$ws.Range("J73").Value = 1600
$ws.Range("K73").Value = 4800
$ws.Range("L73").Value = 4800
$ws.Range("M73").Value = -3864
$ws.Range("N73").Value = -6672
$ws.Range("H86").Value = 146851.86
$ws.Range("I86").Value = 502899.5
$ws.Range("J86").Value = 4432.8
$ws.Range("K86").Value = 502899.5
$ws.Range("L86").Value = 4432.8
$ws.Range("M86").Value = -501776.5
$ws.Range("N86").Value = -6678.8
$ws.Range("H89").Value = 146851.86
$ws.Range("I89").Value = 502899.5
$ws.Range("J89").Value = 4432.8
$ws.Range("K89").Value = 2514497.5
$ws.Range("L89").Value = 22164
$ws.Range("M89").Value = -2508881.5
$ws.Range("N89").Value = -33396
$ws.Range("H113").Value = 4755
$ws.Range("I113").Value = 3999.25
$ws.Range("J113").Value = 5762.6665
$ws.Range("K113").Value = 3999.25
$ws.Range("L113").Value = 5762.6665
$ws.Range("M113").Value = -745.25
$ws.Range("H115").Value = 1375.5555
$ws.Range("I115").Value = 957.5
$ws.Range("J115").Value = 1495
$ws.Range("K115").Value = 2872.5
$ws.Range("L115").Value = 4485
$ws.Range("M115").Value = -1305.5
$ws.Range("H116").Value = 1019334.7
$ws.Range("I116").Value = 1019334.7
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 1019334.7
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = -1015892.7
$ws.Range("N116").ClearContents()
$ws.Range("H118").Value = 1141
$ws.Range("I118").Value = 833.3333
$ws.Range("J118").Value = 1195.2941
$ws.Range("K118").Value = 2499.9999
$ws.Range("L118").Value = 3585.8823
$ws.Range("M118").Value = -842.9998999999998
$ws.Range("N118").Value = -6899.8823
$ws.Range("H132").Value = 5024.3335
$ws.Range("I132").Value = 4351.56
$ws.Range("J132").Value = 13434
$ws.Range("K132").Value = 13054.68
$ws.Range("L132").Value = 40302
$ws.Range("M132").Value = -10524.68
$ws.Range("H135").Value = 42974
$ws.Range("I135").Value = 3603.842
$ws.Range("J135").Value = 167646.17
$ws.Range("K135").Value = 32434.578
$ws.Range("L135").Value = 1508815.53
$ws.Range("M135").Value = -29899.578
$ws.Range("H137").Value = 3397.5
$ws.Range("I137").Value = 2864.125
$ws.Range("J137").Value = 3930.875
$ws.Range("K137").Value = 8592.375
$ws.Range("L137").Value = 11792.625
$ws.Range("M137").Value = -6042.375
$ws.Range("H138").Value = 6935.905
$ws.Range("I138").Value = 12256.75
$ws.Range("J138").Value = 6375.816
$ws.Range("K138").Value = 36770.25
$ws.Range("L138").Value = 19127.448
$ws.Range("M138").Value = -31630.25
$ws.Range("N138").Value = -29407.448
$ws.Range("H141").Value = 2699.4
$ws.Range("I141").Value = 2766.0417
$ws.Range("J141").Value = 1100
$ws.Range("K141").Value = 8298.125100000001
$ws.Range("L141").Value = 3300
$ws.Range("M141").Value = -3118.125100000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 500
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 500
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 500
$ws.Range("N4").Value = -732
$ws.Range("H61").Value = 6214.2856
$ws.Range("I61").Value = 5566.8
$ws.Range("J61").Value = 7833
$ws.Range("K61").Value = 5566.8
$ws.Range("L61").Value = 7833
$ws.Range("M61").Value = -5354.8
$ws.Range("N61").Value = -8257
$ws.Range("H63").Value = 2795.4167
$ws.Range("I63").Value = 2654.5
$ws.Range("J63").Value = 3500
$ws.Range("K63").Value = 2654.5
$ws.Range("L63").Value = 3500
$ws.Range("M63").Value = -1968.5
$ws.Range("H66").Value = 2795.4167
$ws.Range("I66").Value = 2654.5
$ws.Range("J66").Value = 3500
$ws.Range("K66").Value = 13272.5
$ws.Range("L66").Value = 17500
$ws.Range("M66").Value = -9840.5
$ws.Range("H74").Value = 39910.11
$ws.Range("I74").Value = 46163.434
$ws.Range("J74").Value = 3953.5
$ws.Range("K74").Value = 46163.434
$ws.Range("L74").Value = 3953.5
$ws.Range("M74").Value = -45289.434
$ws.Range("H77").Value = 39910.11
$ws.Range("I77").Value = 46163.434
$ws.Range("J77").Value = 3953.5
$ws.Range("K77").Value = 230817.17
$ws.Range("L77").Value = 19767.5
$ws.Range("M77").Value = -226449.17
$ws.Range("H88").Value = 1931.5
$ws.Range("I88").Value = 1300
$ws.Range("J88").Value = 2057.8
$ws.Range("K88").Value = 1300
$ws.Range("L88").Value = 2057.8
$ws.Range("M88").Value = -894
$ws.Range("N88").Value = -2869.8
$ws.Range("H91").Value = 1931.5
$ws.Range("I91").Value = 1300
$ws.Range("J91").Value = 2057.8
$ws.Range("K91").Value = 1300
$ws.Range("L91").Value = 2057.8
$ws.Range("M91").Value = 104
$ws.Range("N91").Value = -4865.8
$ws.Range("H132").Value = 18320.549
$ws.Range("I132").Value = 20226.836
$ws.Range("J132").Value = 3342.5715
$ws.Range("K132").Value = 60680.508
$ws.Range("L132").Value = 10027.7145
$ws.Range("M132").Value = -58150.508
$ws.Range("N132").Value = -15087.7145
$ws.Range("H136").Value = 6214.2856
$ws.Range("I136").Value = 5566.8
$ws.Range("J136").Value = 7833
$ws.Range("K136").Value = 16700.4
$ws.Range("L136").Value = 23499
$ws.Range("M136").Value = -14150.4
$ws.Range("N136").Value = -28599
$ws.Range("H139").Value = 69062.664
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 69062.664
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 69062.664
$ws.Range("N139").Value = -79342.664

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3152
$ws.Range("I134").Value = 2493.2
$ws.Range("J134").Value = 4250
$ws.Range("K134").Value = 7479.599999999999
$ws.Range("L134").Value = 12750
$ws.Range("M134").Value = -4944.599999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2424.4736
$ws.Range("I31").Value = 2104
$ws.Range("J31").Value = 5148.5
$ws.Range("K31").Value = 2104
$ws.Range("L31").Value = 5148.5
$ws.Range("M31").Value = -1809
$ws.Range("H34").Value = 2424.4736
$ws.Range("I34").Value = 2104
$ws.Range("J34").Value = 5148.5
$ws.Range("K34").Value = 2104
$ws.Range("L34").Value = 5148.5
$ws.Range("M34").Value = -1902
$ws.Range("H107").Value = 107428.89
$ws.Range("I107").Value = 115667.695
$ws.Range("J107").Value = 324.5
$ws.Range("K107").Value = 115667.695
$ws.Range("L107").Value = 324.5
$ws.Range("M107").Value = -113747.695

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 1418600
$ws.Range("I129").Value = 3166.6667
$ws.Range("J129").Value = 2834033.2
$ws.Range("K129").Value = 9500.000100000001
$ws.Range("L129").Value = 8502099.600000001
$ws.Range("M129").Value = -4500.000100000001
$ws.Range("N129").Value = -8512099.600000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 5374.75
$ws.Range("I80").Value = 3499
$ws.Range("J80").Value = 6000
$ws.Range("K80").Value = 3499
$ws.Range("L80").Value = 6000
$ws.Range("M80").Value = -2501
$ws.Range("H83").Value = 5374.75
$ws.Range("I83").Value = 3499
$ws.Range("J83").Value = 6000
$ws.Range("K83").Value = 17495
$ws.Range("L83").Value = 30000
$ws.Range("M83").Value = -12503
$ws.Range("H132").Value = 75291.75999999999
$ws.Range("I132").Value = 44940.848
$ws.Range("J132").Value = 338333
$ws.Range("K132").Value = 134822.544
$ws.Range("L132").Value = 1014999
$ws.Range("M132").Value = -132292.544
$ws.Range("H141").Value = 78000
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 78000
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 78000
$ws.Range("N141").Value = -88360

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3898.1365
$ws.Range("I22").Value = 2293.7144
$ws.Range("J22").Value = 4646.8667
$ws.Range("K22").Value = 2293.7144
$ws.Range("L22").Value = 4646.8667
$ws.Range("M22").Value = -1998.7144
$ws.Range("N22").Value = -5236.8667
$ws.Range("H27").Value = 3898.1365
$ws.Range("I27").Value = 2293.7144
$ws.Range("J27").Value = 4646.8667
$ws.Range("K27").Value = 2293.7144
$ws.Range("L27").Value = 4646.8667
$ws.Range("M27").Value = -2186.7144
$ws.Range("N27").Value = -4860.8667
$ws.Range("H55").Value = 712.4194
$ws.Range("I55").Value = 475.77777
$ws.Range("J55").Value = 1040.0769
$ws.Range("K55").Value = 475.77777
$ws.Range("L55").Value = 1040.0769
$ws.Range("M55").Value = -302.77777
$ws.Range("N55").Value = -1386.0769
$ws.Range("H82").Value = 3451
$ws.Range("I82").Value = 2128.4285
$ws.Range("J82").Value = 4376.8
$ws.Range("K82").Value = 2128.4285
$ws.Range("L82").Value = 4376.8
$ws.Range("M82").Value = -1767.4285
$ws.Range("N82").Value = -5098.8
$ws.Range("H85").Value = 3451
$ws.Range("I85").Value = 2128.4285
$ws.Range("J85").Value = 4376.8
$ws.Range("K85").Value = 2128.4285
$ws.Range("L85").Value = 4376.8
$ws.Range("M85").Value = -880.4285
$ws.Range("N85").Value = -6872.8
$ws.Range("H136").Value = 6096.6
$ws.Range("I136").Value = 5541.5
$ws.Range("J136").Value = 7206.8
$ws.Range("K136").Value = 16624.5
$ws.Range("L136").Value = 21620.4
$ws.Range("M136").Value = -14074.5
$ws.Range("H138").Value = 65776.336
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 65776.336
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 65776.336
$ws.Range("N138").Value = -76056.336

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 4000
$ws.Range("I14").Value = 4000
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 4000
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = -3832
$ws.Range("N14").ClearContents()
$ws.Range("H62").Value = 1111898
$ws.Range("I62").Value = 4500
$ws.Range("J62").Value = 1388747.5
$ws.Range("K62").Value = 4500
$ws.Range("L62").Value = 1388747.5
$ws.Range("M62").Value = -3876
$ws.Range("H65").Value = 1111898
$ws.Range("I65").Value = 4500
$ws.Range("J65").Value = 1388747.5
$ws.Range("K65").Value = 22500
$ws.Range("L65").Value = 6943737.5
$ws.Range("M65").Value = -19380
$ws.Range("H81").Value = 629485.6
$ws.Range("I81").Value = 3982.2307
$ws.Range("J81").Value = 3340000.2
$ws.Range("K81").Value = 7964.4614
$ws.Range("L81").Value = 6680000.4
$ws.Range("M81").Value = -6903.4614
$ws.Range("H84").Value = 629485.6
$ws.Range("I84").Value = 3982.2307
$ws.Range("J84").Value = 3340000.2
$ws.Range("K84").Value = 39822.307
$ws.Range("L84").Value = 33400002
$ws.Range("M84").Value = -34518.307
$ws.Range("H132").Value = 26075.62
$ws.Range("I132").Value = 27412.31
$ws.Range("J132").Value = 20060.5
$ws.Range("K132").Value = 82236.93000000001
$ws.Range("L132").Value = 60181.5
$ws.Range("M132").Value = -79706.93000000001
$ws.Range("H136").Value = 13725249
$ws.Range("I136").Value = 1747676.6
$ws.Range("J136").Value = 41672916
$ws.Range("K136").Value = 5243029.800000001
$ws.Range("L136").Value = 125018748
$ws.Range("M136").Value = -5240479.800000001
$ws.Range("N136").Value = -125023848
